$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenario_info")
$ws.Activate()

# --- Insert new row 7: deadHeadTripAllowedModes (group=global) ---
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value2 = "global"
$ws.Range("B7").Value2 = "deadHeadTripAllowedModes"
$ws.Range("C7").Value2 = "pt"
$ws.Range("D7").Clear()
$ws.Range("E7").Value2 = "Deadhead trips are routed using network links that match one of these allowed modes (list sperated with colon: pt,rail)"

# --- Insert new row 12: capacityFactor (group=global) ---
# (original row 11 "seatDurationThreshold" has, after the first insert above, moved to row 12;
#  inserting here pushes it down to row 13 and creates the new row at 12)
$ws.Rows.Item(12).Insert()
$ws.Range("A12").Value2 = "global"
$ws.Range("B12").Value2 = "capacityFactor"
$ws.Range("C12").Value2 = 1.33
$ws.Range("D12").Clear()
$ws.Range("E12").Value2 = "Adjust the passenger capacity of units to reflect deviations in passenger demand."

# --- Fix up the autofilter range to cover the two extra rows (A1:E24 -> A1:E26) ---
$ws.AutoFilterMode = $false
$ws.Range("A1:E26").AutoFilter()

# --- Fix up the workbook-level _FilterDatabase defined name to match ---
$names = $wb.Names
for ($i = 1; $i -le $names.Count; $i++) {
    $n = $names.Item($i)
    if ($n.Name -eq "scenario_info!_FilterDatabase") {
        $n.RefersTo = "=scenario_info!`$A`$1:`$E`$26"
    }
}

# --- Match the view: zoom reset to 100%, and the freshly inserted capacityFactor row selected ---
$excel.ActiveWindow.Zoom = 100
$ws.Rows.Item(12).EntireRow.Select()

Write-Output "done"
